$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (5) below the existing table, mirroring row 4's
# formatting (date style on col A, text-formatted numeric style on col B).
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 45267

# Column B is formatted as text ("@") but this entry, like B4, is a plain
# numeric literal. Briefly reset to the default style so the value is
# stored as a real number, then restore the text format (Excel keeps an
# already-numeric value as a number when you merely change the display
# format afterwards).
$ws.Range("B5").ClearContents()
$ws.Range("B5").Style = "Normal"
$ws.Range("B5").Value = 1016.751
$ws.Range("B5").NumberFormat = "@"

$ws.Range("C5").Value = 30

# Match the author's final selection.
$ws.Range("F15").Select()
